# Added new DP problems: append a "House Robber" row (row 19) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from existing rows so the new row's cell styles match
# the sheet's established pattern (Name=Neutral, Description/Solution=Normal,
# Link=Hyperlink) without minting duplicate style records.
$ws.Range("A6").Copy()
$ws.Range("A19").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("B18").Copy()
$ws.Range("B19").PasteSpecial(-4122)

$ws.Range("C18").Copy()
$ws.Range("C19").PasteSpecial(-4122)

$ws.Range("D18").Copy()
$ws.Range("D19").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Fill in the new problem's data (Name, Link, Description, Solution order
# matches how the source workbook's shared-string table grew).
$ws.Range("A19").Value = "House Robber"
$ws.Range("D19").Value = "https://leetcode.com/problems/house-robber/"
$ws.Range("B19").Value = "Return maximum ammount without robbing adjacent houses"
$ws.Range("C19").Value = "Use a bottom-up approach using for loop and an array. Subproblem: dp[n] = Math.max(dp[n-1], dp[n-2] + nums[n])"

# Turn the URL text into a real hyperlink, then restore the Hyperlink cell
# style (Hyperlinks.Add applies its own font formatting).
$ws.Hyperlinks.Add($ws.Range("D19"), "https://leetcode.com/problems/house-robber/")
$ws.Range("D19").Style = "Hyperlink"

# Scroll the view down one row, matching the author's saved scroll position.
$win = $excel.ActiveWindow
$win.ScrollRow = 2
$win.ScrollColumn = 1
